$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 and E2 stay numeric - straightforward value assignment.
$ws.Range("B2").Value = 0.5742
$ws.Range("E2").Value = 0.6961000000000001

# C2 and D2 are stored as text (numeric-looking strings), matching the
# original workbook's cell type. Temporarily force a Text number format so
# the COM "smart" input doesn't coerce the string into a Number, then reset
# the style back to Normal so no residual formatting is left on the cells.
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C2").Value = "0.9115"
$ws.Range("D2").Value = "0.5742"
$ws.Range("C2:D2").Style = "Normal"
